# [ttml] Address editorial note for padding on content element example.
#
# Re-apply the AutoFilter on the ed-notes table (Status = "O", Assigned to
# = "nigel") to review the remaining open notes, then close out the
# "enhance padding example" note (row 51, column E = Status) now that the
# second padding example has been added. The filter is applied before the
# status edit so the already-filtered rows' hidden/visible state is not
# recomputed (matches normal Excel behaviour: editing a cell does not
# automatically re-run AutoFilter).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A1:H106")
$rng.AutoFilter(5, @("O"), 7)
$rng.AutoFilter(8, @("nigel"), 7)

$ws.Range("E51").Value = "C"

# Move the selection/view to E63 (one of the rows still visible under the
# filter) which also resets the scrolled "topLeftCell" back to default.
$ws.Range("E63").Select()
